# Auto-generated edit script: update market-price / profit columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match
# a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2634.6875
$ws.Range("I2").Value = 871.7143
$ws.Range("J2").Value = 4005.889
$ws.Range("K2").Value = 871.7143
$ws.Range("L2").Value = 4005.889
$ws.Range("M2").Value = -758.7143
$ws.Range("N2").Value = -4231.889
$ws.Range("H17").Value = 2321.3928
$ws.Range("J17").Value = 2321.3928
$ws.Range("L17").Value = 6964.178400000001
$ws.Range("N17").Value = -7300.178400000001
$ws.Range("H42").Value = 1134.1666
$ws.Range("J42").Value = 2244.3333
$ws.Range("L42").Value = 6732.999899999999
$ws.Range("N42").Value = -7192.999899999999
$ws.Range("H62").Value = 3895.8462
$ws.Range("I62").Value = 3564.7
$ws.Range("K62").Value = 3564.7
$ws.Range("M62").Value = -2940.7
$ws.Range("H65").Value = 3895.8462
$ws.Range("I65").Value = 3564.7
$ws.Range("K65").Value = 17823.5
$ws.Range("M65").Value = -14703.5
$ws.Range("H70").Value = 3771.4707
$ws.Range("I70").Value = 2796.4
$ws.Range("J70").Value = 4177.75
$ws.Range("K70").Value = 8389.200000000001
$ws.Range("L70").Value = 12533.25
$ws.Range("M70").Value = -8119.200000000001
$ws.Range("N70").Value = -13073.25
$ws.Range("H73").Value = 3771.4707
$ws.Range("I73").Value = 2796.4
$ws.Range("J73").Value = 4177.75
$ws.Range("K73").Value = 8389.200000000001
$ws.Range("L73").Value = 12533.25
$ws.Range("M73").Value = -7453.200000000001
$ws.Range("N73").Value = -14405.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1673.3055
$ws.Range("I32").Value = 1556.0282
$ws.Range("K32").Value = 1556.0282
$ws.Range("M32").Value = -1269.0282
$ws.Range("H63").Value = 4131.5
$ws.Range("I63").Value = 3665.3333
$ws.Range("J63").Value = 4411.2
$ws.Range("K63").Value = 3665.3333
$ws.Range("L63").Value = 4411.2
$ws.Range("M63").Value = -2979.3333
$ws.Range("N63").Value = -5783.2
$ws.Range("H66").Value = 4131.5
$ws.Range("I66").Value = 3665.3333
$ws.Range("J66").Value = 4411.2
$ws.Range("K66").Value = 18326.6665
$ws.Range("L66").Value = 22056
$ws.Range("M66").Value = -14894.6665
$ws.Range("N66").Value = -28920
$ws.Range("H97").Value = 467.5
$ws.Range("I97").Value = 439
$ws.Range("K97").Value = 439
$ws.Range("M97").Value = 57
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 125000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2940.111
$ws.Range("I31").Value = 2075
$ws.Range("J31").Value = 5189.4
$ws.Range("K31").Value = 2075
$ws.Range("L31").Value = 5189.4
$ws.Range("M31").Value = -1780
$ws.Range("N31").Value = -5779.4
$ws.Range("H34").Value = 2940.111
$ws.Range("I34").Value = 2075
$ws.Range("J34").Value = 5189.4
$ws.Range("K34").Value = 2075
$ws.Range("L34").Value = 5189.4
$ws.Range("M34").Value = -1873
$ws.Range("N34").Value = -5593.4
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H75").Value = 14500
$ws.Range("J75").Value = 14500
$ws.Range("L75").Value = 14500
$ws.Range("N75").Value = -16496
$ws.Range("H78").Value = 14500
$ws.Range("J78").Value = 14500
$ws.Range("L78").Value = 43500
$ws.Range("N78").Value = -53484
$ws.Range("H99").Value = 2250
$ws.Range("I99").Value = 2250
$ws.Range("K99").Value = 2250
$ws.Range("M99").Value = -752
$ws.Range("H120").Value = 40326
$ws.Range("J120").Value = 40326
$ws.Range("L120").Value = 40326
$ws.Range("N120").Value = -47584
$ws.Range("H126").Value = 2250
$ws.Range("I126").Value = 2250
$ws.Range("K126").Value = 6750
$ws.Range("M126").Value = -4280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1656.5
$ws.Range("J55").Value = 1822.4286
$ws.Range("L55").Value = 5467.2858
$ws.Range("N55").Value = -5821.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2588
$ws.Range("I31").Value = 2588
$ws.Range("K31").Value = 2588
$ws.Range("M31").Value = -2296
$ws.Range("H37").Value = 2588
$ws.Range("I37").Value = 2588
$ws.Range("K37").Value = 2588
$ws.Range("M37").Value = -2311
$ws.Range("H46").Value = 16666.666
$ws.Range("J46").Value = 17272.727
$ws.Range("L46").Value = 17272.727
$ws.Range("N46").Value = -17584.727
$ws.Range("H102").Value = 1813.1111
$ws.Range("I102").Value = 1040.5625
$ws.Range("J102").Value = 7993.5
$ws.Range("K102").Value = 1040.5625
$ws.Range("L102").Value = 7993.5
$ws.Range("M102").Value = 581.4375
$ws.Range("N102").Value = -11237.5
$ws.Range("H126").Value = 5024.375
$ws.Range("I126").Value = 6332.5
$ws.Range("K126").Value = 18997.5
$ws.Range("M126").Value = -16527.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 30241
$ws.Range("J6").Value = 30241
$ws.Range("L6").Value = 30241
$ws.Range("N6").Value = -30465
$ws.Range("H30").Value = 3553
$ws.Range("I30").Value = 260
$ws.Range("J30").Value = 20018
$ws.Range("K30").Value = 260
$ws.Range("L30").Value = 20018
$ws.Range("M30").Value = -152
$ws.Range("N30").Value = -20234
$ws.Range("H35").Value = 1372
$ws.Range("I35").Value = 1372
$ws.Range("K35").Value = 1372
$ws.Range("M35").Value = -1036
$ws.Range("H46").Value = 3457.7896
$ws.Range("I46").Value = 2440
$ws.Range("J46").Value = 3821.2856
$ws.Range("K46").Value = 2440
$ws.Range("L46").Value = 3821.2856
$ws.Range("M46").Value = -2252
$ws.Range("N46").Value = -4197.2856
$ws.Range("H55").Value = 200.125
$ws.Range("J55").Value = 202.3077
$ws.Range("L55").Value = 202.3077
$ws.Range("N55").Value = -548.3077000000001
$ws.Range("H61").Value = 7500
$ws.Range("I61").Value = 4997.5
$ws.Range("K61").Value = 4997.5
$ws.Range("M61").Value = -4795.5
$ws.Range("H113").Value = 7500
$ws.Range("I113").Value = 4997.5
$ws.Range("K113").Value = 4997.5
$ws.Range("M113").Value = -2827.5
$ws.Range("H122").Value = 3704.2
$ws.Range("J122").Value = 4505
$ws.Range("L122").Value = 13515
$ws.Range("N122").Value = -18415
$ws.Range("H127").Value = 76331.664
$ws.Range("J127").Value = 76331.664
$ws.Range("L127").Value = 76331.664
$ws.Range("N127").Value = -86251.664
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3025.8333
$ws.Range("I136").Value = 2830.5
$ws.Range("K136").Value = 8491.5
$ws.Range("M136").Value = -5941.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9710
$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55630
$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -57184
$ws.Range("H81").Value = 5111.1113
$ws.Range("I81").Value = 5111.1113
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10222.2226
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -9161.222599999999
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 5111.1113
$ws.Range("I84").Value = 5111.1113
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 51111.113
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -45807.113
$ws.Range("N84").ClearContents()
$ws.Range("H96").Value = 2361.625
$ws.Range("I96").Value = 1890
$ws.Range("J96").Value = 3147.6667
$ws.Range("K96").Value = 1890
$ws.Range("L96").Value = 3147.6667
$ws.Range("M96").Value = -517
$ws.Range("N96").Value = -5893.6667
$ws.Range("H107").Value = 192
$ws.Range("I107").Value = 192
$ws.Range("K107").Value = 576
$ws.Range("M107").Value = 1344
$ws.Range("H122").Value = 3839.75
$ws.Range("I122").Value = 3068.5715
$ws.Range("K122").Value = 9205.7145
$ws.Range("M122").Value = -6755.7145
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
